$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header "Searching Time (s)" -> "Search Time (s)" (column D, header row)
$ws.Range("D1").Value = "Search Time (s)"

# 2. Rename "DIAMOND" -> "* DIAMOND" (A5)
$ws.Range("A5").Value = "* DIAMOND"

# 3. Rename "MMseqs2" -> "* MMseqs2" (A6)
$ws.Range("A6").Value = "* MMseqs2"

# 4. Row heights: row 1 and row 6 go from 15 to 13.8
$ws.Range("A1:G1").RowHeight = 13.8
$ws.Range("A6:G6").RowHeight = 13.8

# 5. Update the active selection to A6
$ws.Range("A6").Select()
